$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 539.8182
$ws.Range("I19").Value = 410.5
$ws.Range("J19").Value = 613.7143
$ws.Range("K19").Value = 410.5
$ws.Range("L19").Value = 613.7143
$ws.Range("M19").Value = -235.5
$ws.Range("N19").Value = -963.7143

$ws.Range("H33").Value = 215.64285
$ws.Range("I33").Value = 193.63637
$ws.Range("J33").Value = 296.33334
$ws.Range("K33").Value = 193.63637
$ws.Range("L33").Value = 296.33334
$ws.Range("M33").Value = 35.36363
$ws.Range("N33").Value = -754.33334

$ws.Range("H38").Value = 1170.0605
$ws.Range("I38").Value = 148.36363
$ws.Range("J38").Value = 1680.909
$ws.Range("K38").Value = 445.09089
$ws.Range("L38").Value = 5042.727000000001
$ws.Range("M38").Value = -73.09089
$ws.Range("N38").Value = -5786.727000000001

$ws.Range("H39").Value = 437.75
$ws.Range("I39").Value = 136.72728
$ws.Range("J39").Value = 1100
$ws.Range("K39").Value = 410.18184
$ws.Range("L39").Value = 3300
$ws.Range("M39").Value = -114.18184
$ws.Range("N39").Value = -3892

$ws.Range("H82").Value = 750
$ws.Range("I82").Value = 750
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2250
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1844

$ws.Range("H85").Value = 750
$ws.Range("I85").Value = 750
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2250
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -846

$ws.Range("H113").Value = 4204436.5
$ws.Range("I113").Value = 11907420
$ws.Range("J113").Value = 2809.0908
$ws.Range("K113").Value = 11907420
$ws.Range("L113").Value = 2809.0908
$ws.Range("M113").Value = -11904166
$ws.Range("N113").Value = -9317.0908

$ws.Range("H125").Value = 1930.1818
$ws.Range("I125").Value = 1032
$ws.Range("J125").Value = 2020
$ws.Range("K125").Value = 9288
$ws.Range("L125").Value = 18180
$ws.Range("M125").Value = -6828
$ws.Range("N125").Value = -23100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2261
$ws.Range("I3").Value = 505
$ws.Range("J3").Value = 2700
$ws.Range("K3").Value = 505
$ws.Range("L3").Value = 2700
$ws.Range("M3").Value = -390
$ws.Range("N3").Value = -2930

$ws.Range("H8").Value = 5006
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 5006
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 5006
$ws.Range("N8").Value = -5294

$ws.Range("H10").Value = 4000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 4000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 4000
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -4340

$ws.Range("H12").Value = 950
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 1500
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = -227
$ws.Range("N12").Value = -1846

$ws.Range("H13").Value = 4004
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 4004
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 4004
$ws.Range("N13").Value = -4292

$ws.Range("H14").Value = 1936.4
$ws.Range("I14").Value = 452.5
$ws.Range("J14").Value = 2925.6667
$ws.Range("K14").Value = 452.5
$ws.Range("L14").Value = 2925.6667
$ws.Range("M14").Value = -277.5
$ws.Range("N14").Value = -3275.6667

$ws.Range("H16").Value = 1800
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1800
$ws.Range("N16").Value = -2374

$ws.Range("H19").Value = 3004.5
$ws.Range("I19").Value = 3000
$ws.Range("J19").Value = 3009
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3009
$ws.Range("M19").Value = -2771
$ws.Range("N19").Value = -3467

$ws.Range("H21").Value = 2677.3333
$ws.Range("I21").Value = 3007.5
$ws.Range("J21").Value = 2017
$ws.Range("K21").Value = 3007.5
$ws.Range("L21").Value = 2017
$ws.Range("M21").Value = -2633.5
$ws.Range("N21").Value = -2765

$ws.Range("H27").Value = 5179
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 5179
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 5179
$ws.Range("N27").Value = -5547

$ws.Range("H30").Value = 6560
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 7450
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 7450
$ws.Range("M30").Value = -2850
$ws.Range("N30").Value = -7750

$ws.Range("H74").Value = 3199.8572
$ws.Range("I74").Value = 4800
$ws.Range("J74").Value = 2933.1667
$ws.Range("K74").Value = 4800
$ws.Range("L74").Value = 2933.1667
$ws.Range("M74").Value = -3926
$ws.Range("N74").Value = -4681.1667

$ws.Range("H77").Value = 3199.8572
$ws.Range("I77").Value = 4800
$ws.Range("J77").Value = 2933.1667
$ws.Range("K77").Value = 24000
$ws.Range("L77").Value = 14665.8335
$ws.Range("M77").Value = -19632
$ws.Range("N77").Value = -23401.8335

$ws.Range("H122").Value = 1329.7037
$ws.Range("I122").Value = 1384.25
$ws.Range("J122").Value = 893.3333
$ws.Range("K122").Value = 4152.75
$ws.Range("L122").Value = 2679.9999
$ws.Range("M122").Value = -1702.75
$ws.Range("N122").Value = -7579.9999

$ws.Range("H129").Value = 49856.715
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49856.715
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49856.715
$ws.Range("N129").Value = -59856.715

$ws.Range("H132").Value = 892577.8
$ws.Range("I132").Value = 1046.1637
$ws.Range("J132").Value = 5350236
$ws.Range("K132").Value = 3138.4911
$ws.Range("L132").Value = 16050708
$ws.Range("M132").Value = -608.4911000000002
$ws.Range("N132").Value = -16055768

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 579
$ws.Range("I10").Value = 165
$ws.Range("J10").Value = 993
$ws.Range("K10").Value = 165
$ws.Range("L10").Value = 993
$ws.Range("M10").Value = -25
$ws.Range("N10").Value = -1273

$ws.Range("H12").Value = 321.6
$ws.Range("I12").Value = 103.333336
$ws.Range("J12").Value = 649
$ws.Range("K12").Value = 103.333336
$ws.Range("L12").Value = 649
$ws.Range("M12").Value = 64.666664
$ws.Range("N12").Value = -985

$ws.Range("H20").Value = 1213.9565
$ws.Range("I20").Value = 1331.7037
$ws.Range("J20").Value = 1046.6316
$ws.Range("K20").Value = 1331.7037
$ws.Range("L20").Value = 1046.6316
$ws.Range("M20").Value = -1084.7037
$ws.Range("N20").Value = -1540.6316

$ws.Range("H134").Value = 1661688
$ws.Range("I134").Value = 1608.7142
$ws.Range("J134").Value = 10113001
$ws.Range("K134").Value = 4826.142599999999
$ws.Range("L134").Value = 30339003
$ws.Range("M134").Value = -2291.142599999999
$ws.Range("N134").Value = -30344073

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 12347920
$ws.Range("I132").Value = 1985.1666
$ws.Range("J132").Value = 22224668
$ws.Range("K132").Value = 5955.4998
$ws.Range("L132").Value = 66674004
$ws.Range("M132").Value = -3425.4998
$ws.Range("N132").Value = -66679064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 13.8125
$ws.Range("I2").Value = 15.923077
$ws.Range("J2").Value = 4.6666665
$ws.Range("K2").Value = 95.538462
$ws.Range("L2").Value = 27.999999
$ws.Range("M2").Value = 17.461538
$ws.Range("N2").Value = -253.999999

$ws.Range("H4").Value = 200260
$ws.Range("I4").Value = 333400
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 1000200
$ws.Range("L4").Value = 1650
$ws.Range("M4").Value = -1000088
$ws.Range("N4").Value = -1874

$ws.Range("H7").Value = 78.85714
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 130.66667
$ws.Range("K7").Value = 120
$ws.Range("L7").Value = 392.00001
$ws.Range("M7").Value = -8
$ws.Range("N7").Value = -616.00001

$ws.Range("H11").Value = 159.71428
$ws.Range("I11").Value = 54.5
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 163.5
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = -23.5
$ws.Range("N11").Value = -1180

$ws.Range("H12").Value = 27517.95
$ws.Range("I12").Value = 4.6666665
$ws.Range("J12").Value = 32373.234
$ws.Range("K12").Value = 13.9999995
$ws.Range("L12").Value = 97119.702
$ws.Range("M12").Value = 159.0000005
$ws.Range("N12").Value = -97465.702

$ws.Range("H26").Value = 220
$ws.Range("I26").Value = 100
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = 300
$ws.Range("L26").Value = 1200
$ws.Range("M26").Value = -12
$ws.Range("N26").Value = -1776

$ws.Range("H68").Value = 3081.681
$ws.Range("I68").Value = 736.88464
$ws.Range("J68").Value = 5984.7617
$ws.Range("K68").Value = 2210.65392
$ws.Range("L68").Value = 17954.2851
$ws.Range("M68").Value = -1399.65392
$ws.Range("N68").Value = -19576.2851

$ws.Range("H71").Value = 3081.681
$ws.Range("I71").Value = 736.88464
$ws.Range("J71").Value = 5984.7617
$ws.Range("K71").Value = 6631.96176
$ws.Range("L71").Value = 53862.8553
$ws.Range("M71").Value = -2575.96176
$ws.Range("N71").Value = -61974.8553

$ws.Range("H107").Value = 860.8099999999999
$ws.Range("I107").Value = 179.31818
$ws.Range("J107").Value = 1053.0256
$ws.Range("K107").Value = 537.9545400000001
$ws.Range("L107").Value = 3159.0768
$ws.Range("M107").Value = 1382.04546
$ws.Range("N107").Value = -6999.0768

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 483.6
$ws.Range("I3").Value = 103.5
$ws.Range("J3").Value = 2004
$ws.Range("K3").Value = 103.5
$ws.Range("L3").Value = 2004
$ws.Range("M3").Value = 12.5
$ws.Range("N3").Value = -2236

$ws.Range("H13").Value = 365
$ws.Range("I13").Value = 62.5
$ws.Range("J13").Value = 566.6667
$ws.Range("K13").Value = 62.5
$ws.Range("L13").Value = 566.6667
$ws.Range("M13").Value = 76.5
$ws.Range("N13").Value = -844.6667

$ws.Range("H70").Value = 6529.8203
$ws.Range("I70").Value = 7546.5713
$ws.Range("J70").Value = 3941.7273
$ws.Range("K70").Value = 7546.5713
$ws.Range("L70").Value = 3941.7273
$ws.Range("M70").Value = -7276.5713
$ws.Range("N70").Value = -4481.7273

$ws.Range("H73").Value = 6529.8203
$ws.Range("I73").Value = 7546.5713
$ws.Range("J73").Value = 3941.7273
$ws.Range("K73").Value = 7546.5713
$ws.Range("L73").Value = 3941.7273
$ws.Range("M73").Value = -6610.5713
$ws.Range("N73").Value = -5813.7273

$ws.Range("H122").Value = 21610008
$ws.Range("I122").Value = 33763868
$ws.Range("J122").Value = 3143.4443
$ws.Range("K122").Value = 101291604
$ws.Range("L122").Value = 9430.332900000001
$ws.Range("M122").Value = -101289154
$ws.Range("N122").Value = -14330.3329

$ws.Range("H132").Value = 5346.927
$ws.Range("I132").Value = 2418.5454
$ws.Range("J132").Value = 17426.5
$ws.Range("K132").Value = 7255.6362
$ws.Range("L132").Value = 52279.5
$ws.Range("M132").Value = -4725.6362
$ws.Range("N132").Value = -57339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8823.786
$ws.Range("I122").Value = 10709.818
$ws.Range("J122").Value = 1908.3334
$ws.Range("K122").Value = 32129.454
$ws.Range("L122").Value = 5725.0002
$ws.Range("M122").Value = -29679.454
$ws.Range("N122").Value = -10625.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9800
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 9800
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 9800
$ws.Range("N15").Value = -10376

$ws.Range("H70").Value = 23915
$ws.Range("I70").Value = 23980
$ws.Range("J70").Value = 23909.092
$ws.Range("K70").Value = 23980
$ws.Range("L70").Value = 23909.092
$ws.Range("M70").Value = -23665
$ws.Range("N70").Value = -24539.092

$ws.Range("H73").Value = 23915
$ws.Range("I73").Value = 23980
$ws.Range("J73").Value = 23909.092
$ws.Range("K73").Value = 23980
$ws.Range("L73").Value = 23909.092
$ws.Range("M73").Value = -22888
$ws.Range("N73").Value = -26093.092

$ws.Range("H122").Value = 75385.71000000001
$ws.Range("I122").Value = 500000
$ws.Range("J122").Value = 4616.6665
$ws.Range("K122").Value = 1500000
$ws.Range("L122").Value = 13849.9995
$ws.Range("M122").Value = -1497550
$ws.Range("N122").Value = -18749.9995
